$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = "IMAGE_RETRIEVAL_BY_METADATA"

# Row 17
$ws.Range("C17").Value = "IMAGE_RETRIEVAL_BY_METADATA"

# Row 21
$ws.Range("D21").Value = $true

# Row 25
$ws.Range("C25").Value = "IMAGE_RETRIEVAL_BY_CAPTION"

# Row 34
$ws.Range("C34").Value = "GEOGRAPHY_QA"
$ws.Range("D34").Value = $true

# Row 46
$ws.Range("C46").Value = "GEOGRAPHY_QA"
$ws.Range("D46").Value = $true

# Row 47
$ws.Range("C47").Value = "None"

# Row 48
$ws.Range("C48").Value = "IMAGE_RETRIEVAL_BY_IMAGE"

# Row 49
$ws.Range("C49").Value = "OBJECT_COUNTING"

# Row 51
$ws.Range("C51").Value = "IMAGE_RETRIEVAL_BY_IMAGE"

# Row 53
$ws.Range("C53").Value = "None"

# Row 54
$ws.Range("C54").Value = "OBJECT_COUNTING"

# Row 56
$ws.Range("C56").Value = "IMAGE_RETRIEVAL_BY_IMAGE"

# Row 58
$ws.Range("C58").Value = "IMAGE_SEGMENTATION"

# Row 59
$ws.Range("C59").Value = "None"

# Row 60
$ws.Range("C60").Value = "IMAGE_RETRIEVAL_BY_IMAGE"

# Row 61
$ws.Range("C61").Value = "OBJECT_COUNTING"
